$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 337, shifting rows 337:346 down to 338:347.
$ws.Rows.Item(337).Insert()

# Populate the new row 337 - it's a copy of the (now shifted) row 338 data
# (same market/region/variety/quality) but with an updated date and prices.
$ws.Range("A337").Value = 11
$ws.Range("B337").Value = "Vega Monumental Concepción"
$ws.Range("C337").Value = "Bíobío"
$ws.Range("D337").Value = 44747
$ws.Range("E337").Value = 8
$ws.Range("F337").Value = 100112006
$ws.Range("G337").Value = "Repollo"
$ws.Range("H337").Value = "Copenhague"
$ws.Range("I337").Value = "Primera"
$ws.Range("J337").Value = 1100
$ws.Range("K337").Value = 1200
$ws.Range("L337").Value = 1300
$ws.Range("M337").Value = 1255
$ws.Range("N337").Value = "$/unidad"
$ws.Range("O337").Value = "Región Metropolitana"
$ws.Range("P337").Value = 1255
$ws.Range("Q337").Value = 1
$ws.Range("R337").Value = "Hortaliza"

# Match style/number format of the date column used elsewhere (column D).
$ws.Range("D337").NumberFormat = $ws.Range("D338").NumberFormat
